$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("C","D","E","F","G","H","I","K","L","N","O")

$data = @{
  2 = @(0.2415770933482833,0.02587126749661905,0.1437606120375037,0.6290542452549559,0.4688686173606129,0.6138120906210389,0.6608872145715878,1.062424887895702,0.1671444806557574,1.056361425878436,2.123412239915922)
  3 = @(0.2355361060694605,0.02380182233109451,0.1399977085992141,0.6268990016914273,0.4682065611508719,0.6174449417600698,0.6601296682002413,0.9374696694845284,0.1626114239204028,1.049576209229173,2.129334901191299)
  4 = @(0.231932996382767,0.02252125527817128,0.1377613719837107,0.626012177325272,0.4681723109783817,0.6199830171014824,0.6600854742578903,0.8605628739133522,0.1599234841592363,1.045776958195177,2.134394430886545)
  5 = @(0.2304914901315129,0.02199695043540117,0.1368687225242802,0.6257604898610651,0.4682518077922708,0.6210946380124938,0.6601732802718772,0.829178909986382,0.1588521325239896,1.044321434083969,2.13681375185385)
  6 = @(0.2302537505712507,0.02190974214516217,0.1367216279900489,0.6257253211887388,0.4682706480967695,0.6212838930475755,0.6601942503912923,0.8239650505950067,0.1586756859546981,1.044085357469228,2.137237062852023)
  7 = @(0.231913447099231,0.02251419425312662,0.1377492577264263,0.6260083388890791,0.4681730049113142,0.6199976956676352,0.6600862300437171,0.860139793050962,0.159908938308412,1.045756952552992,2.134425611506302)
  8 = @(0.2394722106394056,0.02515979953057723,0.1424478072267092,0.62822048123288,0.4685629855412898,0.6150008858313214,0.6605386233109911,1.019379870365441,0.1655616993867142,1.053945991985586,2.125158845845561)
  9 = @(0.2551325357524945,0.03026798509608852,0.1522483714222602,0.6360259750588071,0.47228911440925,0.6076416555410447,0.6647683627657486,1.330100497428703,0.177403287428362,1.072898433004099,2.118293558126595)
  10 = @(0.2671447448893502,0.03397110810480086,0.1598057945356572,0.6438819065664418,0.4768438791421801,0.6037220689096614,0.6699188443031261,1.55734304083262,0.1865653653608348,1.088566395020933,2.120168215883922)
  11 = @(0.27271863043552,0.03564469520084401,0.1633212692099946,0.6479180583428032,0.4793131463221698,0.6022619391563637,0.6727067252050105,1.660475090163061,0.1908340136603925,1.096068707174624,2.122529069941692)
  12 = @(0.2748449604692098,0.03627683252943115,0.1646636119605986,0.649513041473007,0.4803055070994304,0.6017554592553722,0.6738264653914996,1.699491609146435,0.1924649248555426,1.098963166003983,2.123640318703735)
  13 = @(0.2743863244575095,0.0361407627690653,0.1643740210845763,0.6491665708888732,0.4800892331309257,0.6018624732659248,0.6735824606996417,1.691090406000001,0.1921130359300065,1.098337419726477,2.123391323263689)
  14 = @(0.2728932525228345,0.03569673406042284,0.1634314821778915,0.6480479437465334,0.4793936390125708,0.6022193400011986,0.6727975633383636,1.663685766020365,0.1909678999121809,1.096305766233726,2.122616136937182)
  15 = @(0.2719807335053019,0.03542454246683491,0.1628555952349018,0.6473714253694638,0.4789750352029642,0.6024439791426346,0.6723251312457776,1.646894691932175,0.190268354301196,1.095068274490799,2.122169616010837)
  16 = @(0.2667826708391772,0.03386151130834492,0.1595776069671402,0.643627447126029,0.4766905159804367,0.6038239900664735,0.6697456069135086,1.550598033180052,0.1862884248994305,1.08808360949223,2.120044307755535)
  17 = @(0.2636217797604559,0.03289980492747446,0.1575865016945386,0.6414491397896853,0.4753909153034783,0.6047532909669115,0.6682771408512096,1.491459395584741,0.1838726603692322,1.083894464364704,2.119127027463122)
  18 = @(0.2618140300439507,0.03234562522722939,0.1564485740966646,0.6402397582156567,0.4746808030938183,0.6053181937419794,0.6674743910402867,1.457421817618808,0.1824926632321962,1.081520288267356,2.118741355181982)
  19 = @(0.2612037317054501,0.0321578133734377,0.1560645471111926,0.6398377550692942,0.4744467866592714,0.6055146800573254,0.6672097838197217,1.445893476423976,0.1820270500006558,1.080722510089018,2.11863513675948)
  20 = @(0.263957195717694,0.03300228723349363,0.1577977026909991,0.6416765192493159,0.4755253896716027,0.6046512196479625,0.6684291274742407,1.497757162083758,0.1841288409623445,1.08433675462399,2.119209982193269)
  21 = @(0.2733313809745255,0.0358272001282316,0.1637080275728238,0.6483747040578223,0.4795963953018116,0.6021132591909151,0.6730263685337761,1.67173620735241,0.1913038616925888,1.096901063486399,2.122837928929897)
  22 = @(0.2795489449918875,0.03766402161954829,0.1676354959853228,0.6531404706999524,0.482591117670907,0.6007252430520964,0.6764041725751042,1.785222670350606,0.1960774640440519,1.105424195603263,2.126475446602058)
  23 = @(0.2762222284175664,0.03668454898418361,0.1655334243631472,0.6505613525848446,0.4809621520203962,0.601441283640753,0.6745672058365102,1.724673692658314,0.1935219962755923,1.100846859437297,2.124418028211352)
  24 = @(0.2638055245981832,0.03295595897293424,0.1577021975822248,0.6415735871959001,0.475464478428762,0.6046972706638485,0.6683602850847663,1.49491005997794,0.1840129941233641,1.084136688575057,2.119172037049964)
  25 = @(0.2508067500459674,0.0288947558685777,0.149534321284122,0.633542397703458,0.4709628330583939,0.6093713516338113,0.6632657220701645,1.246218474773343,0.1741187412480087,1.067463646127976,2.1189376672923)
}

foreach ($r in $data.Keys) {
  $vals = $data[$r]
  for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range("$($cols[$i])$r").Value = $vals[$i]
  }
}

Write-Output "done"